$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.901.04'
$ws.Range('E2').Value = '  -2.32%  '
$ws.Range('D3').Value = '2.488.31'
$ws.Range('E3').Value = '  -5.30%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''577.90'
$ws.Range('E5').Value = '  -2.79%  '
$ws.Range('D6').Value = '''169.50'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.517'
$ws.Range('E8').Value = '  -3.19%  '
$ws.Range('D9').Value = '2.488.40'
$ws.Range('E9').Value = '  -5.30%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '''0.347'
$ws.Range('E12').Value = '  -4.85%  '
$ws.Range('D13').Value = '''5.07'
$ws.Range('E13').Value = '  -3.03%  '
$ws.Range('D14').Value = '''26.27'
$ws.Range('D15').Value = '2.936.55'
$ws.Range('E15').Value = '  -5.47%  '
$ws.Range('E16').Value = '  -5.14%  '
$ws.Range('D17').Value = '65.717.76'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').Value = '2.508.75'
$ws.Range('E18').Value = '  -4.48%  '
$ws.Range('D19').Value = '''11.11'
$ws.Range('E19').Value = '  -7.76%  '
$ws.Range('D20').Value = '''7.56'
$ws.Range('E20').Value = '  -5.31%  '
$ws.Range('D21').Value = '''342.70'
$ws.Range('E21').Value = '  -3.82%  '
$ws.Range('D22').Value = '''4.16'
$ws.Range('E22').Value = '  -3.56%  '
$ws.Range('E23').Value = '  -2.96%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '''1.94'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '''68.72'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('D27').Value = '''9.88'
$ws.Range('E27').Value = '  -3.72%  '
$ws.Range('D28').Value = '''0.998'
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D29').Value = '2.613.65'
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('D30').Value = '0.0₃0957'
$ws.Range('E30').Value = '  -4.87%  '
$ws.Range('D31').Value = '''520.18'
$ws.Range('E31').Value = '  -4.48%  '
$ws.Range('D32').Value = '''8.02'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('E33').Value = '  -3.23%  '
$ws.Range('D34').Value = '''1.81'
$ws.Range('E34').Value = '  -4.86%  '
$ws.Range('D35').Value = '''0.131'
$ws.Range('E35').Value = '  -3.78%  '
$ws.Range('D36').Value = '''0.999'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').Value = '''157.61'
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  -4.01%  '
$ws.Range('D39').Value = '''18.41'
$ws.Range('E39').Value = '  -3.22%  '
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('E41').Value = '  -4.05%  '
$ws.Range('E42').Value = '  -3.27%  '
$ws.Range('D43').Value = '''5.03'
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = '''145.66'
$ws.Range('E46').Value = '  -4.63%  '
$ws.Range('D47').Value = '''0.551'
$ws.Range('E47').Value = '  -5.07%  '
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('D49').Value = '0.0₆0269'
$ws.Range('E49').Value = '  -9.41%  '
$ws.Range('D50').Value = '''1.69'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').Value = '''0.0746'
$ws.Range('E51').Value = '  -2.98%  '
